# Auto-generated edit script applying numeric corrections to the
# Leve profit-calculation columns (H,I,J,K,L,M,N) across all 8 job sheets,
# per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ==== Sheet: ALC ====
# Row 28
$ws_ALC.Range("H28").Value = 1918.1818
$ws_ALC.Range("J28").Value = 1184.5
$ws_ALC.Range("L28").Value = 1184.5
$ws_ALC.Range("N28").Value = -2154.5
# Row 53
$ws_ALC.Range("H53").Value = 431
$ws_ALC.Range("I53").Value = 336.68182
$ws_ALC.Range("J53").Value = 727.4286
$ws_ALC.Range("K53").Value = 336.68182
$ws_ALC.Range("L53").Value = 727.4286
$ws_ALC.Range("M53").Value = 300.31818
$ws_ALC.Range("N53").Value = -2001.4286
# Row 62
$ws_ALC.Range("H62").Value = 22228488
$ws_ALC.Range("I62").Value = 44451530
$ws_ALC.Range("K62").Value = 44451530
$ws_ALC.Range("M62").Value = -44450906
# Row 65
$ws_ALC.Range("H65").Value = 22228488
$ws_ALC.Range("I65").Value = 44451530
$ws_ALC.Range("K65").Value = 222257650
$ws_ALC.Range("M65").Value = -222254530
# Row 76
$ws_ALC.Range("H76").Value = 7150160
$ws_ALC.Range("I76").Value = 12505871
$ws_ALC.Range("K76").Value = 12505871
$ws_ALC.Range("M76").Value = -12505556
# Row 79
$ws_ALC.Range("H79").Value = 7150160
$ws_ALC.Range("I79").Value = 12505871
$ws_ALC.Range("K79").Value = 12505871
$ws_ALC.Range("M79").Value = -12504779
# Row 107
$ws_ALC.Range("H107").Value = 18519392
$ws_ALC.Range("I107").Value = 22728244
$ws_ALC.Range("K107").Value = 22728244
$ws_ALC.Range("M107").Value = -22726324
# Row 111
$ws_ALC.Range("H111").Value = 4993.6
$ws_ALC.Range("I111").Value = 4992.8887
$ws_ALC.Range("J111").Value = 5000
$ws_ALC.Range("K111").Value = 14978.6661
$ws_ALC.Range("L111").Value = 15000
$ws_ALC.Range("M111").Value = -11911.6661
$ws_ALC.Range("N111").Value = -21134
# Row 118
$ws_ALC.Range("H118").Value = 276701
$ws_ALC.Range("J118").Value = 999
$ws_ALC.Range("L118").Value = 2997
$ws_ALC.Range("N118").Value = -6311
# Row 132
$ws_ALC.Range("H132").Value = 1288.4419
$ws_ALC.Range("I132").Value = 992.7805
$ws_ALC.Range("K132").Value = 2978.3415
$ws_ALC.Range("M132").Value = -448.3415

# ==== Sheet: ARM ====
# Row 32
$ws_ARM.Range("H32").Value = 31241.652
$ws_ARM.Range("I32").Value = 32312.309
$ws_ARM.Range("K32").Value = 32312.309
$ws_ARM.Range("M32").Value = -32025.309
# Row 63
$ws_ARM.Range("H63").Value = 7491.1724
$ws_ARM.Range("I63").Value = 2783.4
$ws_ARM.Range("K63").Value = 2783.4
$ws_ARM.Range("M63").Value = -2097.4
# Row 66
$ws_ARM.Range("H66").Value = 7491.1724
$ws_ARM.Range("I66").Value = 2783.4
$ws_ARM.Range("K66").Value = 13917
$ws_ARM.Range("M66").Value = -10485
# Row 122
$ws_ARM.Range("H122").Value = 1656.1765
$ws_ARM.Range("I122").Value = 1414.4546
$ws_ARM.Range("J122").Value = 2099.3333
$ws_ARM.Range("K122").Value = 4243.3638
$ws_ARM.Range("L122").Value = 6297.999899999999
$ws_ARM.Range("M122").Value = -1793.3638
$ws_ARM.Range("N122").Value = -11197.9999
# Row 132
$ws_ARM.Range("H132").Value = 3178359.2
$ws_ARM.Range("I132").Value = 3848836.8
$ws_ARM.Range("K132").Value = 11546510.4
$ws_ARM.Range("M132").Value = -11543980.4

# ==== Sheet: BSM ====
# Row 20
$ws_BSM.Range("H20").Value = 2281.1428
$ws_BSM.Range("I20").Value = 1889.6666
$ws_BSM.Range("K20").Value = 1889.6666
$ws_BSM.Range("M20").Value = -1642.6666
# Row 63
$ws_BSM.Range("H63").Value = 0
$ws_BSM.Range("J63").Value = 0
$ws_BSM.Range("L63").Value = 0
$ws_BSM.Range("N63").ClearContents()
# Row 66
$ws_BSM.Range("H66").Value = 0
$ws_BSM.Range("J66").Value = 0
$ws_BSM.Range("L66").Value = 0
$ws_BSM.Range("N66").ClearContents()
# Row 92
$ws_BSM.Range("H92").Value = 0
$ws_BSM.Range("J92").Value = 0
$ws_BSM.Range("L92").Value = 0
$ws_BSM.Range("N92").ClearContents()
# Row 107
$ws_BSM.Range("H107").Value = 1760.579
$ws_BSM.Range("I107").Value = 1350.3077
$ws_BSM.Range("J107").Value = 2649.5
$ws_BSM.Range("K107").Value = 1350.3077
$ws_BSM.Range("L107").Value = 2649.5
$ws_BSM.Range("M107").Value = 569.6922999999999
$ws_BSM.Range("N107").Value = -6489.5
# Row 131
$ws_BSM.Range("H131").Value = 57329.332
$ws_BSM.Range("J131").Value = 57329.332
$ws_BSM.Range("L131").Value = 57329.332
$ws_BSM.Range("N131").Value = -67409.33199999999

# ==== Sheet: CRP ====
# Row 22
$ws_CRP.Range("H22").Value = 432.35
$ws_CRP.Range("I22").Value = 369.25
$ws_CRP.Range("K22").Value = 369.25
$ws_CRP.Range("M22").Value = -19.25
# Row 31
$ws_CRP.Range("H31").Value = 3892.2092
$ws_CRP.Range("I31").Value = 2168.2173
$ws_CRP.Range("J31").Value = 5874.8
$ws_CRP.Range("K31").Value = 2168.2173
$ws_CRP.Range("L31").Value = 5874.8
$ws_CRP.Range("M31").Value = -1873.2173
$ws_CRP.Range("N31").Value = -6464.8
# Row 34
$ws_CRP.Range("H34").Value = 3892.2092
$ws_CRP.Range("I34").Value = 2168.2173
$ws_CRP.Range("J34").Value = 5874.8
$ws_CRP.Range("K34").Value = 2168.2173
$ws_CRP.Range("L34").Value = 5874.8
$ws_CRP.Range("M34").Value = -1966.2173
$ws_CRP.Range("N34").Value = -6278.8
# Row 62
$ws_CRP.Range("H62").Value = 9024
$ws_CRP.Range("I62").Value = 6795.75
$ws_CRP.Range("J62").Value = 11252.25
$ws_CRP.Range("K62").Value = 6795.75
$ws_CRP.Range("L62").Value = 11252.25
$ws_CRP.Range("M62").Value = -6171.75
$ws_CRP.Range("N62").Value = -12500.25
# Row 65
$ws_CRP.Range("H65").Value = 9024
$ws_CRP.Range("I65").Value = 6795.75
$ws_CRP.Range("J65").Value = 11252.25
$ws_CRP.Range("K65").Value = 33978.75
$ws_CRP.Range("L65").Value = 56261.25
$ws_CRP.Range("M65").Value = -30858.75
$ws_CRP.Range("N65").Value = -62501.25
# Row 107
$ws_CRP.Range("H107").Value = 100000620
$ws_CRP.Range("I107").Value = 111111660
$ws_CRP.Range("J107").Value = 1300
$ws_CRP.Range("K107").Value = 111111660
$ws_CRP.Range("L107").Value = 1300
$ws_CRP.Range("M107").Value = -111109740
$ws_CRP.Range("N107").Value = -5140

# ==== Sheet: CUL ====
# Row 3
$ws_CUL.Range("H3").Value = 8017.875
$ws_CUL.Range("I3").Value = 7020.5713
$ws_CUL.Range("J3").Value = 14999
$ws_CUL.Range("K3").Value = 21061.7139
$ws_CUL.Range("L3").Value = 44997
$ws_CUL.Range("M3").Value = -20949.7139
$ws_CUL.Range("N3").Value = -45221
# Row 22
$ws_CUL.Range("H22").Value = 903.2222
$ws_CUL.Range("J22").Value = 1159.8
$ws_CUL.Range("L22").Value = 3479.4
$ws_CUL.Range("N22").Value = -3817.4
# Row 27
$ws_CUL.Range("H27").Value = 903.2222
$ws_CUL.Range("J27").Value = 1159.8
$ws_CUL.Range("L27").Value = 3479.4
$ws_CUL.Range("N27").Value = -3683.4
# Row 131
$ws_CUL.Range("H131").Value = 11115189
$ws_CUL.Range("J131").Value = 4726
$ws_CUL.Range("L131").Value = 14178
$ws_CUL.Range("N131").Value = -24258
# Row 133
$ws_CUL.Range("H133").Value = 4749.5
$ws_CUL.Range("I133").Value = 4749.5
$ws_CUL.Range("K133").Value = 14248.5
$ws_CUL.Range("M133").Value = -9188.5
# Row 134
$ws_CUL.Range("H134").Value = 3337.8
$ws_CUL.Range("I134").Value = 3337.8
$ws_CUL.Range("K134").Value = 10013.4
$ws_CUL.Range("M134").Value = -4943.400000000001

# ==== Sheet: GSM ====
# Row 107
$ws_GSM.Range("H107").Value = 679.4138
$ws_GSM.Range("J107").Value = 326.5
$ws_GSM.Range("L107").Value = 326.5
$ws_GSM.Range("N107").Value = -4166.5
# Row 113
$ws_GSM.Range("H113").Value = 2365.16
$ws_GSM.Range("I113").Value = 2227.3157
$ws_GSM.Range("J113").Value = 2801.6667
$ws_GSM.Range("K113").Value = 2227.3157
$ws_GSM.Range("L113").Value = 2801.6667
$ws_GSM.Range("M113").Value = -57.31570000000011
$ws_GSM.Range("N113").Value = -7141.6667
# Row 122
$ws_GSM.Range("H122").Value = 3447.3572
$ws_GSM.Range("I122").Value = 1905.6538
$ws_GSM.Range("J122").Value = 5952.625
$ws_GSM.Range("K122").Value = 5716.9614
$ws_GSM.Range("L122").Value = 17857.875
$ws_GSM.Range("M122").Value = -3266.9614
$ws_GSM.Range("N122").Value = -22757.875
# Row 132
$ws_GSM.Range("H132").Value = 6714.4585
$ws_GSM.Range("I132").Value = 4249.294
$ws_GSM.Range("K132").Value = 12747.882
$ws_GSM.Range("M132").Value = -10217.882

# ==== Sheet: LTW ====
# Row 61
$ws_LTW.Range("H61").Value = 1942.88
$ws_LTW.Range("I61").Value = 2183.6191
$ws_LTW.Range("K61").Value = 2183.6191
$ws_LTW.Range("M61").Value = -1981.6191
# Row 113
$ws_LTW.Range("H113").Value = 1942.88
$ws_LTW.Range("I113").Value = 2183.6191
$ws_LTW.Range("K113").Value = 2183.6191
$ws_LTW.Range("M113").Value = -13.61909999999989

# ==== Sheet: WVR ====
# Row 54
$ws_WVR.Range("H54").Value = 203332.67
$ws_WVR.Range("J54").Value = 204999
$ws_WVR.Range("L54").Value = 204999
$ws_WVR.Range("N54").Value = -206039
# Row 113
$ws_WVR.Range("H113").Value = 731.11536
$ws_WVR.Range("I113").Value = 929
$ws_WVR.Range("K113").Value = 2787
$ws_WVR.Range("M113").Value = -617
# Row 136
$ws_WVR.Range("H136").Value = 4083636.5
$ws_WVR.Range("I136").Value = 4927286.5
$ws_WVR.Range("J136").Value = 5994
$ws_WVR.Range("K136").Value = 14781859.5
$ws_WVR.Range("L136").Value = 17982
$ws_WVR.Range("M136").Value = -14779309.5
$ws_WVR.Range("N136").Value = -23082

